# The flashcards sheet is being reset: all question/answer rows plus the
# bold/bordered/centered header style applied to row 1 are removed, leaving
# a single blank sheet (dimension collapses back down to A1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear every cell (values + the header's bold font / thin border / center
# alignment formatting) so the sheet goes back to being completely empty.
$ws.Cells.Clear()
